$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 05:29"

# --- Update COVID figures for countries whose row stayed in place ---

# Row 33: Kazajistan
$ws.Range("B33").Value = 106498
$ws.Range("C33").Value = 73
$ws.Range("E33").Value = 4971
$ws.Range("H33").Value = 1634

# Row 39: Belgica
$ws.Range("B39").Value = 89141
$ws.Range("C39").Value = 372
$ws.Range("D39").Value = 18602
$ws.Range("E39").Value = 60627
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 9912

# Row 75: Australia
$ws.Range("B75").Value = 26465
$ws.Range("C75").Value = 91
$ws.Range("D75").Value = 22725
$ws.Range("E75").Value = 2959
$ws.Range("G75").Value = 11
$ws.Range("H75").Value = 781

# Row 172: Islas Turcas y Caicos
$ws.Range("B172").Value = 614
$ws.Range("C172").Value = 15
$ws.Range("E172").Value = 345

# Row 173: San Martin (Parte Holandesa)
$ws.Range("B173").Value = 527
$ws.Range("C173").Value = 11
$ws.Range("E173").Value = 187

# --- Insert Birmania as a new ranked row right after Yemen (row 149), ---
# --- pushing Nueva Zelanda, Georgia and Uruguay down one row each.     ---
# Guyana (row 154) and everything below stays where it is.

# Row 153 now becomes Uruguay, taking the old row 152 data (Uruguay)
$ws.Range("A153").Value = "Uruguay"
$ws.Range("B153").Value = 1712
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 1476
$ws.Range("E153").Value = 191
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 45

# Row 152 now becomes Georgia, taking the old row 151 data (Georgia)
$ws.Range("A152").Value = "Georgia"
$ws.Range("B152").Value = 1729
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 1321
$ws.Range("E152").Value = 389
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 19

# Row 151 now becomes Nueva Zelanda, taking the old row 150 data (Nueva Zelanda)
$ws.Range("A151").Value = "Nueva Zelanda"
$ws.Range("B151").Value = 1788
$ws.Range("C151").Value = 6
$ws.Range("D151").Value = 1639
$ws.Range("E151").Value = 125
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 24

# Row 150 becomes Birmania with its own new figures
$ws.Range("A150").Value = "Birmania"
$ws.Range("B150").Value = 1807
$ws.Range("C150").Value = 98
$ws.Range("D150").Value = 460
$ws.Range("E150").Value = 1335
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 2
$ws.Range("H150").Value = 12
